$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.523.47'
$ws.Range('E2').Value = '  +0.40%  '
$ws.Range('D3').Value = '1.848.55'
$ws.Range('E3').Value = '  +0.33%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '263.53'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.48%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5225'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.23%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3236'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.66%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06766'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.18%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.68'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.99%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.7706'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07768'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.12%  '
$ws.Range('D13').Value = '1.859.33'
$ws.Range('E13').Value = '  +0.84%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '88.32'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.27%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.011'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.14%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.001'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.08%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '13.93'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.94%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.001'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.07%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007931'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.47%  '
$ws.Range('D20').Value = '26.557.82'
$ws.Range('E20').Value = '  +0.39%  '
$ws.Range('D21').Value = '2.091.31'
$ws.Range('E21').Value = '  +0.32%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.616'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.89%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.425'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.10%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.961'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.12%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '142.78'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.48%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.179'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -6.53%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.677'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.83%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '16.98'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.11%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '111.64'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.61%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.165'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.20%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08740'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.39%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.110'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.17%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04820'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.01%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.128'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.38%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.875'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.06%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7150'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +4.38%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.104'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.71%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01783'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.07%  '
$ws.Range('E39').Value = '  -1.14%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.4840'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.37%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '112.33'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.93%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8963'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.75%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.046'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.12%  '
$ws.Range('E44').Value = '  -0.01%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '7.615'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.02%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4163'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.94%  '
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05901'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.15%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.030'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.58%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '34.91'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.14%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.1227'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.39%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.8841'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.00%  '
